# Actualización automática hashcode
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B9" = "3c173f6c33219d20ff80d6484c04951a"
    "B15" = "71c45d49f862b1f742de0ed7f20f578e"
    "B17" = "eb742a8b7a274769449d10cc70c50362"
    "B24" = "a641771cb9ce92af6589fa4466ffb1f6"
    "B79" = "12f2372106f41e89a1b7421a7f2ebce5"
    "B133" = "bb31d92226927ee3d096f20d51505076"
    "B136" = "1726c225f6c647908c3111e0f76c434f"
    "B150" = "5fc04dfc411f3aee10b1f4527a2eb929"
    "B162" = "6a38fe89afbde4c1e16a9178a905d1bb"
    "B180" = "b8abe84582fe1dd8e80cb732ec3dbda0"
    "B183" = "1ab3dfe9d627d6391bbdb545c98f6600"
    "B191" = "76258f163628e4ef93faa285cb7e34f7"
    "B198" = "332b4d635c3adc2008f68f66c4bf8fbc"
    "B200" = "122fbf96a0d2014a0c43ef15b9afabf3"
    "B227" = "f843e7138a69cf36ea4a98900026d31f"
    "B232" = "c8cd1474f29108d901faf76c814618c3"
    "B293" = "b33b7c774497249d0d3e3a8a69119531"
    "B302" = "b90b0a084a07d7e058c1ad06b189406b"
    "B339" = "dd7f2c72288ad89e02c68d5064673d6d"
    "B415" = "e76114d8232bc912a4cc66a523a2a7b6"
    "B460" = "cb3839402c4dbcb396059113ff253b37"
    "B478" = "643efad5a5e7ac50a38d12919d4aeacb"
    "B480" = "560b9902838cade8ba14010b7573f9b1"
    "B500" = "2bc00e21214a125e2f24299f389417ea"
    "B501" = "0d7b477812b30e50c64e93ffae6dbbcc"
    "B502" = "d3d143d051c1500894bec1918b6360c5"
    "B503" = "15b74f135bb25aa373625983be9c33c6"
    "B504" = "c600014ecbeff376e2860bffe842415a"
    "B506" = "dd6565526824a99002fb01872fc6e124"
    "B515" = "4ec0dd30e65481d2465ec25e9b46fd79"
    "B517" = "6ef9a29dbe581b14cf019305a840aa49"
    "B547" = "6cafa3c8c1e7a757077370f6c39320ed"
    "B550" = "c664e2e989756550f718c4e96f3130dc"
    "B563" = "d04c87c8bce96f3f6000d4fdc1f80468"
    "B616" = "26d3ad6dba27e3ebe34758f7188b79a5"
    "B627" = "30d9f61399cea4aab4f0212cb42d8164"
    "B665" = "5fa589b8940045feb5f2bfcb3ef11078"
    "B685" = "6d49c5c66b55659c67ce568f5079f223"
    "B700" = "04bf18e52f0d3df8ff37eaffdb581a7b"
    "B703" = "47b16e19009dc2c81e46ff534969b397"
    "B704" = "c5041d579b96618864c3ba6dce926b4f"
    "B742" = "84d7d00287f261c2a6707731f0a04c6b"
    "B795" = "09b2547196d057257fa8d355bc56555a"
    "B819" = "23f4a822775b8c7f10fa2c24c447b7fc"
    "B830" = "24647f8535dea3acac025fbcb3f286e5"
    "B835" = "7c14b075f4ec9477861832201411932d"
    "B854" = "174523586ba67d40e3592fbbf1678536"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
